# Auto-generated edit script applying market-price data refresh
# to the Raiden_Profits workbook, per sheet/row/column.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2468.739
$ws.Range("I28").Value = 1593.3334
$ws.Range("J28").Value = 3423.7273
$ws.Range("K28").Value = 1593.3334
$ws.Range("L28").Value = 3423.7273
$ws.Range("M28").Value = -1108.3334
$ws.Range("N28").Value = -4393.7273

$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H43").Value = 2692.889
$ws.Range("J43").Value = 3084
$ws.Range("L43").Value = 3084
$ws.Range("N43").Value = -3222

$ws.Range("H62").Value = 5333.3335
$ws.Range("I62").Value = 6000
$ws.Range("K62").Value = 6000
$ws.Range("M62").Value = -5376

$ws.Range("H65").Value = 5333.3335
$ws.Range("I65").Value = 6000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26880

$ws.Range("I74").Value = 4996.25
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4996.25
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -4060.25
$ws.Range("N74").Value = -6872

$ws.Range("H76").Value = 5998
$ws.Range("I76").Value = 6997
$ws.Range("K76").Value = 6997
$ws.Range("M76").Value = -6682

$ws.Range("I77").Value = 4996.25
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 24981.25
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -20301.25
$ws.Range("N77").Value = -34360

$ws.Range("H79").Value = 5998
$ws.Range("I79").Value = 6997
$ws.Range("K79").Value = 6997
$ws.Range("M79").Value = -5905

$ws.Range("H113").Value = 9333.733
$ws.Range("I113").Value = 8855
$ws.Range("J113").Value = 9880.857
$ws.Range("K113").Value = 8855
$ws.Range("L113").Value = 9880.857
$ws.Range("M113").Value = -5601
$ws.Range("N113").Value = -16388.857

$ws.Range("H125").Value = 736.9231
$ws.Range("I125").Value = 226
$ws.Range("J125").Value = 1056.25
$ws.Range("K125").Value = 2034
$ws.Range("L125").Value = 9506.25
$ws.Range("M125").Value = 426
$ws.Range("N125").Value = -14426.25

$ws.Range("H135").Value = 8400.857
$ws.Range("I135").Value = 1557
$ws.Range("K135").Value = 14013
$ws.Range("M135").Value = -11478

$ws.Range("H138").Value = 3291.8494
$ws.Range("I138").Value = 3070.8462
$ws.Range("J138").Value = 3572.1462
$ws.Range("K138").Value = 9212.5386
$ws.Range("L138").Value = 10716.4386
$ws.Range("M138").Value = -4072.5386
$ws.Range("N138").Value = -20996.4386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 614.58826
$ws.Range("I5").Value = 491.30768
$ws.Range("J5").Value = 1015.25
$ws.Range("K5").Value = 491.30768
$ws.Range("L5").Value = 1015.25
$ws.Range("M5").Value = -379.30768
$ws.Range("N5").Value = -1239.25

$ws.Range("H63").Value = 2575.5557
$ws.Range("I63").Value = 2497.3684
$ws.Range("K63").Value = 2497.3684
$ws.Range("M63").Value = -1811.3684

$ws.Range("H66").Value = 2575.5557
$ws.Range("I66").Value = 2497.3684
$ws.Range("K66").Value = 12486.842
$ws.Range("M66").Value = -9054.841999999999

$ws.Range("H97").Value = 1342.15
$ws.Range("J97").Value = 2573
$ws.Range("L97").Value = 2573
$ws.Range("N97").Value = -3565

$ws.Range("H102").Value = 1783.0476
$ws.Range("I102").Value = 822.8125
$ws.Range("K102").Value = 822.8125
$ws.Range("M102").Value = 799.1875

$ws.Range("H122").Value = 2705.6584
$ws.Range("I122").Value = 2635.2856
$ws.Range("K122").Value = 7905.8568
$ws.Range("M122").Value = -5455.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 614.58826
$ws.Range("I4").Value = 491.30768
$ws.Range("J4").Value = 1015.25
$ws.Range("K4").Value = 491.30768
$ws.Range("L4").Value = 1015.25
$ws.Range("M4").Value = -376.30768
$ws.Range("N4").Value = -1245.25

$ws.Range("H20").Value = 999.5
$ws.Range("J20").Value = 999.5
$ws.Range("L20").Value = 999.5
$ws.Range("N20").Value = -1493.5

$ws.Range("H86").Value = 4815.2856
$ws.Range("I86").Value = 4453
$ws.Range("J86").Value = 4960.2
$ws.Range("K86").Value = 4453
$ws.Range("L86").Value = 4960.2
$ws.Range("M86").Value = -3330
$ws.Range("N86").Value = -7206.2

$ws.Range("H89").Value = 4815.2856
$ws.Range("I89").Value = 4453
$ws.Range("J89").Value = 4960.2
$ws.Range("K89").Value = 22265
$ws.Range("L89").Value = 24801
$ws.Range("M89").Value = -16649
$ws.Range("N89").Value = -36033

$ws.Range("H94").Value = 1608.3448
$ws.Range("I94").Value = 633.43475
$ws.Range("K94").Value = 633.43475
$ws.Range("M94").Value = -182.43475

$ws.Range("H105").Value = 3799.2
$ws.Range("I105").Value = 3250
$ws.Range("J105").Value = 4165.3335
$ws.Range("K105").Value = 3250
$ws.Range("L105").Value = 4165.3335
$ws.Range("M105").Value = -1503
$ws.Range("N105").Value = -7659.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 146.83333
$ws.Range("J7").Value = 100
$ws.Range("L7").Value = 100
$ws.Range("N7").Value = -326

$ws.Range("H17").Value = 2999
$ws.Range("I17").Value = 2999
$ws.Range("K17").Value = 2999
$ws.Range("M17").Value = -2825

$ws.Range("H99").Value = 1783.8572
$ws.Range("I99").Value = 1783.8572
$ws.Range("K99").Value = 1783.8572
$ws.Range("M99").Value = -285.8571999999999

$ws.Range("H126").Value = 1783.8572
$ws.Range("I126").Value = 1783.8572
$ws.Range("K126").Value = 5351.571599999999
$ws.Range("M126").Value = -2881.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 161.1
$ws.Range("I2").Value = 211.57143
$ws.Range("J2").Value = 43.333332
$ws.Range("K2").Value = 211.57143
$ws.Range("L2").Value = 43.333332
$ws.Range("M2").Value = -98.57142999999999
$ws.Range("N2").Value = -269.333332

$ws.Range("H70").Value = 7842.55
$ws.Range("I70").Value = 6524.75
$ws.Range("J70").Value = 8721.083000000001
$ws.Range("K70").Value = 6524.75
$ws.Range("L70").Value = 8721.083000000001
$ws.Range("M70").Value = -6254.75
$ws.Range("N70").Value = -9261.083000000001

$ws.Range("H73").Value = 7842.55
$ws.Range("I73").Value = 6524.75
$ws.Range("J73").Value = 8721.083000000001
$ws.Range("K73").Value = 6524.75
$ws.Range("L73").Value = 8721.083000000001
$ws.Range("M73").Value = -5588.75
$ws.Range("N73").Value = -10593.083

$ws.Range("H80").Value = 13298.23
$ws.Range("I80").Value = 5485.375
$ws.Range("K80").Value = 5485.375
$ws.Range("M80").Value = -4487.375

$ws.Range("H83").Value = 13298.23
$ws.Range("I83").Value = 5485.375
$ws.Range("K83").Value = 27426.875
$ws.Range("M83").Value = -22434.875

$ws.Range("H97").Value = 746.53845
$ws.Range("I97").Value = 516.7273
$ws.Range("K97").Value = 516.7273
$ws.Range("M97").Value = -20.72730000000001

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H107").Value = 603.7368
$ws.Range("I107").Value = 690
$ws.Range("J107").Value = 580.73334
$ws.Range("K107").Value = 690
$ws.Range("L107").Value = 580.73334
$ws.Range("M107").Value = 1230
$ws.Range("N107").Value = -4420.73334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2516.8845
$ws.Range("I40").Value = 2512.55
$ws.Range("J40").Value = 2531.3333
$ws.Range("K40").Value = 2512.55
$ws.Range("L40").Value = 2531.3333
$ws.Range("M40").Value = -2376.55
$ws.Range("N40").Value = -2803.3333

$ws.Range("H122").Value = 4247.7026
$ws.Range("I122").Value = 4357.857
$ws.Range("K122").Value = 13073.571
$ws.Range("M122").Value = -10623.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 21150
$ws.Range("I29").Value = 21150
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 21150
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -20860
$ws.Range("N29").ClearContents()
